$wb = $excel.ActiveWorkbook

# --- Status text update: "Ready for handoff" -> "Handed back: in sync with en-US" ---
# This shared string is used on: Overview!E2:F2/E3:F3, and Status column (C2/C3) on
# both the zh-cn and de-de sheets.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

$mdName1 = "0f3b7744-d886-4321-a041-d05818ea3ab5.md"
$mdName2 = "a202f2c0-a649-4e3a-b0a2-dfe85867b23b.md"
$url1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3b4a57034a4476e6e725164f927b894e914c9020/e2e/0f3b7744-d886-4321-a041-d05818ea3ab5.md"
$url2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3b4a57034a4476e6e725164f927b894e914c9020/e2e/a202f2c0-a649-4e3a-b0a2-dfe85867b23b.md"

foreach ($pair in @(
        @{ Name = "zh-cn"; Xlf1 = "0f3b7744-d886-4321-a041-d05818ea3ab5.af497dc83affc8b0c5ccad918c766214d9286b10.zh-cn.xlf"; Xlf2 = "a202f2c0-a649-4e3a-b0a2-dfe85867b23b.d4cf6276af13839f94aab748915850f724005d71.zh-cn.xlf"; HandbackDate = "2016-08-30 00:49:26" },
        @{ Name = "de-de"; Xlf1 = "0f3b7744-d886-4321-a041-d05818ea3ab5.af497dc83affc8b0c5ccad918c766214d9286b10.de-de.xlf"; Xlf2 = "a202f2c0-a649-4e3a-b0a2-dfe85867b23b.d4cf6276af13839f94aab748915850f724005d71.de-de.xlf"; HandbackDate = "2016-08-30 00:49:33" }
    )) {
    $ws = $wb.Worksheets.Item($pair.Name)

    # Status column (C) -> "Handed back: in sync with en-US"
    $ws.Range("C2").Value = "Handed back: in sync with en-US"
    $ws.Range("C3").Value = "Handed back: in sync with en-US"

    # Latest Target File (I) column: now populated with the source .md file name,
    # rendered as a hyperlink pointing at the same target as column A's link.
    $ws.Range("I2").Value = $mdName1
    $ws.Hyperlinks.Add($ws.Range("I2"), $url1, "", "", $mdName1)
    $ws.Range("I2").Font.Underline = 2
    $ws.Range("I2").Font.Color = 15570276

    $ws.Range("I3").Value = $mdName2
    $ws.Hyperlinks.Add($ws.Range("I3"), $url2, "", "", $mdName2)
    $ws.Range("I3").Font.Underline = 2
    $ws.Range("I3").Font.Color = 15570276

    # Latest Handback File (J) column: populated with the generated xlf file name.
    $ws.Range("J2").Value = $pair.Xlf1
    $ws.Range("J3").Value = $pair.Xlf2

    # Latest Handback DateTime (K) column: populated with the handback timestamp.
    $ws.Range("K2").Value = $pair.HandbackDate
    $ws.Range("K3").Value = $pair.HandbackDate

    # Column width updates (I, J widened to fit full file names; Status (C) widened too).
    $ws.Columns.Item(3).ColumnWidth = 29.9777047293527
    $ws.Columns.Item(9).ColumnWidth = 40
    $ws.Columns.Item(10).ColumnWidth = 40
}

# Overview sheet column widths for zh-cn/de-de status columns (E, F).
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

$wb.Save()
